$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style used for plain (non-bold/no-border) data cells, taken from an
# untouched text cell so re-applying it after forcing a cell to text
# format does not leave a stray explicit style behind.
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "35.207.15"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").Value = "1.895.49"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  -0.38%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "245.46"
$cell.Style = $plainStyle
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("E6").Value = "  +7.72%  "

$ws.Range("E7").Value = "  -0.27%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "40.72"
$cell.Style = $plainStyle
$ws.Range("E8").Value = "  -3.61%  "

$ws.Range("E9").Value = "  +2.43%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "53.10"
$cell.Style = $plainStyle
$ws.Range("E10").Value = "  +11.46%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0718"
$cell.Style = $plainStyle
$ws.Range("E11").Value = "  +1.79%  "

$ws.Range("E13").Value = "  -0.62%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "12.57"
$cell.Style = $plainStyle
$ws.Range("E14").Value = "  +0.99%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.701"
$cell.Style = $plainStyle
$ws.Range("E15").Value = "  +1.40%  "

$ws.Range("D16").Value = "1.894.34"
$ws.Range("E16").Value = "  -0.72%  "

$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("D18").Value = "35.184.31"
$ws.Range("E18").Value = "  -0.94%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "71.83"
$cell.Style = $plainStyle
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("E20").Value = "  +0.52%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "239.83"
$cell.Style = $plainStyle
$ws.Range("E21").Value = "  -1.63%  "

$ws.Range("E22").Value = "  +0.79%  "

$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("E24").Value = "  -0.25%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.31"
$cell.Style = $plainStyle
$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("E26").Value = "  +5.27%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "167.63"
$cell.Style = $plainStyle
$ws.Range("E27").Value = "  -2.48%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "8.53"
$cell.Style = $plainStyle
$ws.Range("E28").Value = "  -0.28%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "18.19"
$cell.Style = $plainStyle
$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("E30").Value = "  +3.48%  "

$ws.Range("E31").Value = "  +20.57%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.14"
$cell.Style = $plainStyle
$ws.Range("E32").Value = "  +0.96%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.0565"
$cell.Style = $plainStyle
$ws.Range("E33").Value = "  +0.24%  "

$ws.Range("E34").Value = "  +8.56%  "

$ws.Range("E35").Value = "  -0.31%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.918"
$cell.Style = $plainStyle
$ws.Range("E36").Value = "  -8.10%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "4.08"
$cell.Style = $plainStyle
$ws.Range("E37").Value = "  -1.38%  "

$ws.Range("E38").Value = "  +12.95%  "

$ws.Range("E39").Value = "  -0.73%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "16.45"
$cell.Style = $plainStyle
$ws.Range("E40").Value = "  +6.44%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.09"
$cell.Style = $plainStyle
$ws.Range("E41").Value = "  -1.96%  "

$ws.Range("E42").Value = "  +0.75%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.0632"
$cell.Style = $plainStyle
$ws.Range("E43").Value = "  +6.68%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "89.52"
$cell.Style = $plainStyle
$ws.Range("E44").Value = "  -2.03%  "

$ws.Range("D45").Value = "1.348.35"
$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("E46").Value = "  +2.29%  "

$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("E48").Value = "  +0.87%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "45.68"
$cell.Style = $plainStyle
$ws.Range("E49").Value = "  -14.56%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "12.10"
$cell.Style = $plainStyle
$ws.Range("E50").Value = "  -4.87%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "6.44"
$cell.Style = $plainStyle
$ws.Range("E51").Value = "  -3.16%  "
